# Updates the cryptos price/volume table with the latest scraped snapshot.
# All cells in columns B:E are stored as literal text in the workbook (prices
# such as "1.00" or "38.10" are not numbers), so every write below is forced
# to text (leading apostrophe) and then reset to the default "Normal" style
# so no stray number formatting / quote-prefix styling is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

Set-TextCell 2 4 '64.640.54'
Set-TextCell 2 5 '  +1.07%  '
Set-TextCell 3 4 '3.482.91'
Set-TextCell 3 5 '  +11.78%  '
Set-TextCell 4 4 '1.03'
Set-TextCell 4 5 '  +2.65%  '
Set-TextCell 5 4 '595.57'
Set-TextCell 5 5 '  +1.54%  '
Set-TextCell 6 4 '148.63'
Set-TextCell 6 5 '  +2.03%  '
Set-TextCell 7 4 '1.02'
Set-TextCell 7 5 '  +1.80%  '
Set-TextCell 8 4 '3.184.08'
Set-TextCell 8 5 '  +2.48%  '
Set-TextCell 9 4 '0.536'
Set-TextCell 9 5 '  +1.26%  '
Set-TextCell 10 5 '  +1.39%  '
Set-TextCell 11 4 '6.07'
Set-TextCell 11 5 '  +6.26%  '
Set-TextCell 13 4 '0.0000249'
Set-TextCell 13 5 '  +0.36%  '
Set-TextCell 14 4 '38.10'
Set-TextCell 14 5 '  +2.97%  '
Set-TextCell 15 4 '3.726.38'
Set-TextCell 15 5 '  +2.63%  '
Set-TextCell 17 2 'Polkadot'
Set-TextCell 17 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 17 4 '7.38'
Set-TextCell 17 5 '  +4.18%  '
Set-TextCell 18 2 'WrappedEther'
Set-TextCell 18 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 18 4 '3.233.59'
Set-TextCell 18 5 '  +3.87%  '
Set-TextCell 19 4 '64.668.30'
Set-TextCell 19 5 '  +1.48%  '
Set-TextCell 20 4 '476.54'
Set-TextCell 20 5 '  +2.76%  '
Set-TextCell 21 4 '14.80'
Set-TextCell 21 5 '  +3.79%  '
Set-TextCell 22 4 '0.749'
Set-TextCell 22 5 '  +3.19%  '
Set-TextCell 23 4 '7.76'
Set-TextCell 23 5 '  +4.22%  '
Set-TextCell 24 4 '2.50'
Set-TextCell 24 5 '  +12.18%  '
Set-TextCell 25 4 '13.61'
Set-TextCell 25 5 '  +4.38%  '
Set-TextCell 26 4 '83.24'
Set-TextCell 26 5 '  +2.24%  '
Set-TextCell 27 2 'Dai'
Set-TextCell 27 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 27 4 '1.00'
Set-TextCell 27 5 '  +0.24%  '
Set-TextCell 28 2 'RenderToken'
Set-TextCell 28 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 28 4 '10.02'
Set-TextCell 28 5 '  +7.94%  '
Set-TextCell 29 5 '  +2.43%  '
Set-TextCell 30 4 '2.27'
Set-TextCell 30 5 '  +2.45%  '
Set-TextCell 31 2 'NEARProtocol'
Set-TextCell 31 3 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 31 4 '7.36'
Set-TextCell 31 5 '  +4.82%  '
Set-TextCell 32 2 'FirstDigitalUSD'
Set-TextCell 32 3 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell 32 4 '1.00'
Set-TextCell 32 5 '  +0.02%  '
Set-TextCell 33 5 '  +11.44%  '
Set-TextCell 34 4 '28.47'
Set-TextCell 34 5 '  +6.39%  '
Set-TextCell 35 4 '0.0₃0863'
Set-TextCell 35 5 '  +0.24%  '
Set-TextCell 36 4 '1.08'
Set-TextCell 36 5 '  +4.37%  '
Set-TextCell 37 4 '6.28'
Set-TextCell 37 5 '  +4.55%  '
Set-TextCell 38 4 '2.33'
Set-TextCell 38 5 '  +0.60%  '
Set-TextCell 39 4 '3.33'
Set-TextCell 39 5 '  -4.09%  '
Set-TextCell 40 4 '471.16'
Set-TextCell 40 5 '  +6.68%  '
Set-TextCell 41 4 '51.93'
Set-TextCell 41 5 '  +3.00%  '
Set-TextCell 42 4 '9.44'
Set-TextCell 42 5 '  +8.61%  '
Set-TextCell 43 4 '0.303'
Set-TextCell 43 5 '  +9.56%  '
Set-TextCell 45 4 '2.935.79'
Set-TextCell 45 5 '  +1.39%  '
Set-TextCell 46 5 '  +2.38%  '
Set-TextCell 47 5 '  +7.04%  '
Set-TextCell 48 4 '131.87'
Set-TextCell 48 5 '  +4.30%  '
Set-TextCell 49 5 '  +6.20%  '
Set-TextCell 51 4 '24.98'
Set-TextCell 51 5 '  +2.98%  '
